# fixed some issues with part description and adjusted BOM
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# R1, R6, R7, R73 (row 28): Manufacturer Partno corrected from 0402 to 0603 package
$ws.Range("D28").Value = "RC0603JR-070RL"

# R30 (row 44): Description corrected from "2.49Ohm" to "2.49kOhm"
$ws.Range("C44").Value = "Resistor 2.49kOhm 0402 1%"

# R34 (row 46): Manufacturer Partno corrected to the right part
$ws.Range("D46").Value = "CRGCQ0402F1M0"

# Scroll/selection state as left by the author after editing
$ws.Range("D51").Select()
